$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains its original text formatting
# (values like "1.50" or "0.0220" must not be auto-converted to numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.943.83"
$ws.Range("E2").Value = "  +2.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.433.57"
$ws.Range("E3").Value = "  +4.91%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.63"
$ws.Range("E5").Value = "  +1.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.55"
$ws.Range("E6").Value = "  +5.92%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.431.20"
$ws.Range("E9").Value = "  +4.93%  "

$ws.Range("E10").Value = "  +3.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  +3.87%  "

$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("E13").Value = "  +4.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("E14").Value = "  +10.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.866.54"
$ws.Range("E15").Value = "  +4.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.881.56"
$ws.Range("E16").Value = "  +2.93%  "

$ws.Range("E17").Value = "  +7.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.438.79"
$ws.Range("E18").Value = "  +5.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("E19").Value = "  +5.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.81"
$ws.Range("E20").Value = "  +9.67%  "

$ws.Range("E21").Value = "  +2.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.81"
$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.01"
$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.50"
$ws.Range("E27").Value = "  +11.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.19"
$ws.Range("E28").Value = "  +5.66%  "

$ws.Range("E29").Value = "  +11.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0784"
$ws.Range("E30").Value = "  +7.33%  "

$ws.Range("E31").Value = "  +4.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "171.93"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.31"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.44"
$ws.Range("E34").Value = "  +5.40%  "

$ws.Range("E35").Value = "  +4.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.56"
$ws.Range("E36").Value = "  +4.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.47"
$ws.Range("E37").Value = "  +10.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "364.91"
$ws.Range("E38").Value = "  +14.90%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  +11.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.18"
$ws.Range("E42").Value = "  +3.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.68"
$ws.Range("E43").Value = "  +6.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  +5.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.59"
$ws.Range("E45").Value = "  +8.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0954"
$ws.Range("E46").Value = "  +1.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0515"
$ws.Range("E48").Value = "  +4.67%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.84"
$ws.Range("E49").Value = "  +5.91%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0220"
$ws.Range("E50").Value = "  +3.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0216"
$ws.Range("E51").Value = "  -6.37%  "
